# Rename the two embedded logo pictures that live in the document's
# header/footer stories:
#   - BTec_Logo-Orange picture (in both headers):  image1.jpg -> image2.jpg
#   - PearsonLogo picture      (in both footers):  image2.png -> image1.png
#
# InlineShapes that live inside a Header/Footer story need to be
# re-fetched through $word.Selection after calling .Select() on them --
# renaming the handle obtained directly from
# Sections(n).Headers/Footers(n).Range.InlineShapes(1) works for the
# header story but not for the footer story unless the shape is first
# selected and re-queried via the application Selection object.

$d = $word.ActiveDocument

function Rename-InlineShape($shape, $newName) {
    $shape.Select()
    $selected = $word.Selection.InlineShapes.Item(1)
    $selected.Name = $newName
}

$section = $d.Sections.Item(1)

# --- Headers: BTec_Logo-Orange, image1.jpg -> image2.jpg -------------
Rename-InlineShape $section.Headers.Item(1).Range.InlineShapes.Item(1) "image2.jpg"
Rename-InlineShape $section.Headers.Item(2).Range.InlineShapes.Item(1) "image2.jpg"

# --- Footers: PearsonLogo, image2.png -> image1.png -------------------
Rename-InlineShape $section.Footers.Item(1).Range.InlineShapes.Item(1) "image1.png"
Rename-InlineShape $section.Footers.Item(2).Range.InlineShapes.Item(1) "image1.png"
